$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.221.46"
$ws.Range("E2").Value = "  +5.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.612.23"
$ws.Range("E3").Value = "  +4.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.61"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.73"
$ws.Range("E6").Value = "  +3.12%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.611.64"
$ws.Range("E9").Value = "  +4.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  +13.84%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.349"
$ws.Range("E12").Value = "  +3.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.06"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.114.41"
$ws.Range("E14").Value = "  +5.68%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.84"
$ws.Range("E15").Value = "  +5.10%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000185"
$ws.Range("E16").Value = "  +8.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.233.63"
$ws.Range("E17").Value = "  +5.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.641.83"
$ws.Range("E18").Value = "  +5.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.85"
$ws.Range("E19").Value = "  +5.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.41"
$ws.Range("E20").Value = "  +4.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.96"
$ws.Range("E21").Value = "  +5.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.18"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.45"
$ws.Range("E24").Value = "  +4.12%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("E26").Value = "  +10.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  +7.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.746.59"
$ws.Range("E28").Value = "  +6.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0954"
$ws.Range("E30").Value = "  +6.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "533.66"
$ws.Range("E31").Value = "  +7.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.10"
$ws.Range("E32").Value = "  +4.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").Value = "  +5.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +4.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.41"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.24"
$ws.Range("E38").Value = "  +4.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.96"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.38"
$ws.Range("E40").Value = "  +5.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  +5.13%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("E43").Value = "  +5.25%  "
$ws.Range("E44").Value = "  +5.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.330"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.77"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.63"
$ws.Range("E47").Value = "  +4.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.68"
$ws.Range("E48").Value = "  +3.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0270"
$ws.Range("E49").Value = "  +4.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.532"
$ws.Range("E50").Value = "  +3.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  +5.99%  "
